$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3 data
$ws.Range("A3").Value = 8102000
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1

# Best-fit column A's width now that it holds longer PatternID values (e.g. 8102000)
$ws.Columns.Item(1).ColumnWidth = 8.3

# Update selection to match the diff (C4 is selected next)
$ws.Range("C4").Select()
